$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expediente ")

# 1) Edit existing values on row 2
$ws.Range("C2").Value = "miztli editado 2"
$ws.Range("H2").Value = "no sé"

# 2) Build the new row 3 as a copy of row 2 (same field values for all
#    columns except the email + nombre which differ for the new record),
#    then tweak the two differing cells.
$ws.Range("A2:N2").Copy()
$ws.Range("A3:N3").PasteSpecial(-4104)

$ws.Range("A3").Value = "prueba_carga2@hotmail.com"
$ws.Range("C3").Value = "Citlalli 2"

# 3) Hyperlink the new e-mail cell the same way A2 is hyperlinked, then
#    restore A3's formatting (the hyperlink call stamps a blue/underline
#    style) back to the plain style used by the rest of the data rows.
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:prueba_carga2@hotmail.com", "", "", "prueba_carga2@hotmail.com")
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

$ws.Range("A3").Select() | Out-Null
